$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo in the shared "interface" text used by rows 10 & 11
$ws.Range("G10").Value = "Calibrate, SampleProvider"
$ws.Range("G11").Value = "Calibrate, SampleProvider"

# Row 6: AngleSensor claimed by Lawrie (not finished), Mode=Angle, interface=SamplerProvider
$ws.Range("D6").Value = "Lawrie"
$ws.Range("E6").Value = "N"
$ws.Range("F6").Value = "Angle"
$ws.Range("G6").Value = "SamplerProvider"

# Row 7: BarometricHTSensor claimed by Lawrie (not finished), Mode=Pressure, interface=SampleProvider
$ws.Range("D7").Value = "Lawrie"
$ws.Range("E7").Value = "N"
$ws.Range("F7").Value = "Pressure"
$ws.Range("G7").Value = "SampleProvider"

# Update the active cell selection to reflect where editing left off
$ws.Range("G7").Select()
